$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44232
$ws.Range("J2").Value2 = 250

# Row 4
$ws.Range("D4").Value2 = 44189
$ws.Range("J4").Value2 = 250
$ws.Range("K4").Value2 = 5000
$ws.Range("M4").Value2 = 5500
$ws.Range("O4").Value2 = "Provincia de Quillota"
$ws.Range("P4").Value2 = 344

# Row 5
$ws.Range("D5").Value2 = 44215
$ws.Range("J5").Value2 = 250

# Row 6
$ws.Range("D6").Value2 = 44210
$ws.Range("J6").Value2 = 340

# Row 7
$ws.Range("D7").Value2 = 44186
$ws.Range("J7").Value2 = 160

# Row 8
$ws.Range("D8").Value2 = 44230
$ws.Range("J8").Value2 = 250

# Row 9
$ws.Range("D9").Value2 = 44208
$ws.Range("J9").Value2 = 160

# Row 10
$ws.Range("D10").Value2 = 44204
$ws.Range("J10").Value2 = 430

# Row 11
$ws.Range("D11").Value2 = 44292
$ws.Range("J11").Value2 = 90
$ws.Range("K11").Value2 = 6000
$ws.Range("M11").Value2 = 6000
$ws.Range("O11").Value2 = "Región Metropolitana"
$ws.Range("P11").Value2 = 375

# Row 12
$ws.Range("D12").Value2 = 44187
$ws.Range("J12").Value2 = 160

# Row 13
$ws.Range("D13").Value2 = 44188
$ws.Range("J13").Value2 = 210
